$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SWF")

# --- Add new row 13 data (new "Smith-Woolhiser Overland" example) ---
$ws.Range("A13").Formula = "=A12+1"
$ws.Range("B13").Value = "Smith-Woolhiser Overland"
$ws.Range("C13").Value = 0.000122
$ws.Range("D13").Value = 0.01
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0.0001
$ws.Range("G13").Value = 0.0001
$ws.Range("H13").Value = "CENTIMETERS"
$ws.Range("I13").Value = "MINUTES"

# --- Rebuild the running index formulas for rows 11-12 so the A4:A10 shared series extends through A13 ---
$ws.Range("A11").Formula = "=A10+1"
$ws.Range("A12").Formula = "=A11+1"

# --- Re-lock body rows 2:13 (previously rows 10-12 were left unlocked) ---
$ws.Range("A2:J13").Locked = $true

# --- Selection / window bookkeeping to match the saved state ---
$ws.Range("A1:XFD13").Select()

Write-Host "edit applied"
